$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: politeness_score (B38) was stored as text "4"; store it as a
# true number (4) instead.
$ws.Range("B38").Value = 4

# Row 39 (new): duplicate of row 38's annotation, but for a different
# sentence_purpose (D39 = "APC" instead of "SMY"). politeness_score
# stays a text value "4" here, matching the original row 38 style.
$ws.Range("A39").Value = "Sunsi Wu"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "4"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "well"
$ws.Range("D39").Value = "APC"
$ws.Range("E39").Value = "MET"
$ws.Range("F39").Value = "d0296b92-10f5-497e-8726-aae675ac805b"
$ws.Range("G39").Value = "rJl3yM-Ab_annotated.xlsx"
$ws.Range("H39").Value = "The new method is motivated well and departs from prior work."
